$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.027224318692906
$ws.Cells.Item(2, 4).Value = 1.034947433758511
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.042025972480211
$ws.Cells.Item(2, 9).Value = 1.032428291166399
$ws.Cells.Item(2, 10).Value = 1.032383381785719
$ws.Cells.Item(2, 11).Value = 1.037745445740943
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.044803824631255
$ws.Cells.Item(2, 14).Value = 1.033849484603898
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.028169839486751
$ws.Cells.Item(3, 4).Value = 1.035670738606791
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.042941783818132
$ws.Cells.Item(3, 9).Value = 1.032594693427095
$ws.Cells.Item(3, 10).Value = 1.032969017136981
$ws.Cells.Item(3, 11).Value = 1.038278272787041
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.045530108037946
$ws.Cells.Item(3, 14).Value = 1.034435951624518
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.028781641315332
$ws.Cells.Item(4, 4).Value = 1.036138378863789
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.043534502801539
$ws.Cells.Item(4, 9).Value = 1.032700419575627
$ws.Cells.Item(4, 10).Value = 1.033347351214813
$ws.Cells.Item(4, 11).Value = 1.038621993666104
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.045999561660409
$ws.Cells.Item(4, 14).Value = 1.034814822980136
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.029038838856795
$ws.Cells.Item(5, 4).Value = 1.036334881010815
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.043783711194474
$ws.Cells.Item(5, 9).Value = 1.032744399972488
$ws.Cells.Item(5, 10).Value = 1.033506255693095
$ws.Cells.Item(5, 11).Value = 1.038766240605842
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.046196799056294
$ws.Cells.Item(5, 14).Value = 1.034973953121004
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.029082023202357
$ws.Cells.Item(6, 4).Value = 1.036367869045766
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.043825556063341
$ws.Cells.Item(6, 9).Value = 1.032751757076258
$ws.Cells.Item(6, 10).Value = 1.033532927799192
$ws.Cells.Item(6, 11).Value = 1.038790445417121
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.046229908959826
$ws.Cells.Item(6, 14).Value = 1.035000663104551
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.02878507801793
$ws.Cells.Item(7, 4).Value = 1.036141004905418
$ws.Cells.Item(7, 5).Value = 0.9943035907982488
$ws.Cells.Item(7, 6).Value = 1.043537832622633
$ws.Cells.Item(7, 9).Value = 1.032701009079242
$ws.Cells.Item(7, 10).Value = 1.033349475083588
$ws.Cells.Item(7, 11).Value = 1.038623922097445
$ws.Cells.Item(7, 12).Value = 0.9968970624462087
$ws.Cells.Item(7, 13).Value = 1.046002197632669
$ws.Cells.Item(7, 14).Value = 1.034816949865048
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.02754386418598
$ws.Cells.Item(8, 4).Value = 1.035191957220138
$ws.Cells.Item(8, 5).Value = 0.9929600610674301
$ws.Cells.Item(8, 6).Value = 1.042335447987293
$ws.Cells.Item(8, 9).Value = 1.032484930358237
$ws.Cells.Item(8, 10).Value = 1.032581426025236
$ws.Cells.Item(8, 11).Value = 1.037925734728883
$ws.Cells.Item(8, 12).Value = 0.995817528259106
$ws.Cells.Item(8, 13).Value = 1.045049378422148
$ws.Cells.Item(8, 14).Value = 1.034047810088937
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.025356617191275
$ws.Cells.Item(9, 4).Value = 1.033516712899713
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.04021772678526
$ws.Cells.Item(9, 9).Value = 1.032089289102921
$ws.Cells.Item(9, 10).Value = 1.031223382058165
$ws.Cells.Item(9, 11).Value = 1.036687412816186
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.043366606904942
$ws.Cells.Item(9, 14).Value = 1.032687837543663
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.023898447241012
$ws.Cells.Item(10, 4).Value = 1.032398009288751
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.038806680528128
$ws.Cells.Item(10, 9).Value = 1.031815562060276
$ws.Cells.Item(10, 10).Value = 1.030314946187584
$ws.Cells.Item(10, 11).Value = 1.035856526909318
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.042242281044733
$ws.Cells.Item(10, 14).Value = 1.031778111589946
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.023267050280188
$ws.Cells.Item(11, 4).Value = 1.031913169380913
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.038195878740815
$ws.Cells.Item(11, 9).Value = 1.031694679532572
$ws.Cells.Item(11, 10).Value = 1.029920864270636
$ws.Cells.Item(11, 11).Value = 1.035495489914027
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.04175486085167
$ws.Cells.Item(11, 14).Value = 1.031383470031475
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.023032521925603
$ws.Cells.Item(12, 4).Value = 1.031733014269335
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.037969029506766
$ws.Cells.Item(12, 9).Value = 1.031649424882669
$ws.Cells.Item(12, 10).Value = 1.029774376774617
$ws.Cells.Item(12, 11).Value = 1.035361196604157
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.041573725189382
$ws.Cells.Item(12, 14).Value = 1.031236774506409
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.02308282901828
$ws.Cells.Item(13, 4).Value = 1.031771661028478
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.038017688073614
$ws.Cells.Item(13, 9).Value = 1.031659148145558
$ws.Cells.Item(13, 10).Value = 1.029805803712048
$ws.Cells.Item(13, 11).Value = 1.035390011472421
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.04161258328423
$ws.Cells.Item(13, 14).Value = 1.031268246073696
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.023247664078029
$ws.Cells.Item(14, 4).Value = 1.031898279001502
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.038177126699127
$ws.Cells.Item(14, 9).Value = 1.031690945976859
$ws.Cells.Item(14, 10).Value = 1.029908757765993
$ws.Cells.Item(14, 11).Value = 1.03548439301908
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.041739889872043
$ws.Cells.Item(14, 14).Value = 1.031371346334206
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.023349224493169
$ws.Cells.Item(15, 4).Value = 1.031976284071073
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.038275366073326
$ws.Cells.Item(15, 9).Value = 1.031710490846346
$ws.Cells.Item(15, 10).Value = 1.029972176876959
$ws.Cells.Item(15, 11).Value = 1.035542519707163
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.0418183162828
$ws.Cells.Item(15, 14).Value = 1.031434855507584
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.023940351030315
$ws.Cells.Item(16, 4).Value = 1.032430177478374
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.038847221541405
$ws.Cells.Item(16, 9).Value = 1.031823535029667
$ws.Cells.Item(16, 10).Value = 1.030341084931817
$ws.Cells.Item(16, 11).Value = 1.035880461319224
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.042274617389201
$ws.Cells.Item(16, 14).Value = 1.031804287454194
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.024311149547649
$ws.Cells.Item(17, 4).Value = 1.032714777514523
$ws.Cells.Item(17, 5).Value = 0.989476357848556
$ws.Cells.Item(17, 6).Value = 1.039205983176115
$ws.Cells.Item(17, 9).Value = 1.031893814074367
$ws.Cells.Item(17, 10).Value = 1.030572298103843
$ws.Cells.Item(17, 11).Value = 1.036092107017611
$ws.Cells.Item(17, 12).Value = 0.9930127773699352
$ws.Cells.Item(17, 13).Value = 1.042560688523296
$ws.Cells.Item(17, 14).Value = 1.032035828975446
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.024527430036731
$ws.Cells.Item(18, 4).Value = 1.032880737959202
$ws.Cells.Item(18, 5).Value = 0.9897087662937556
$ws.Cells.Item(18, 6).Value = 1.039415261079052
$ws.Cells.Item(18, 9).Value = 1.031934579204224
$ws.Cells.Item(18, 10).Value = 1.030707091009583
$ws.Cells.Item(18, 11).Value = 1.036215434942782
$ws.Cells.Item(18, 12).Value = 0.9932001317071769
$ws.Cells.Item(18, 13).Value = 1.042727492979607
$ws.Cells.Item(18, 14).Value = 1.032170813302573
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.024601176075318
$ws.Cells.Item(19, 4).Value = 1.032937319031385
$ws.Cells.Item(19, 5).Value = 0.9897880325774034
$ws.Cells.Item(19, 6).Value = 1.039486622507515
$ws.Cells.Item(19, 9).Value = 1.031948440462511
$ws.Cells.Item(19, 10).Value = 1.030753040041268
$ws.Cells.Item(19, 11).Value = 1.036257465954281
$ws.Cells.Item(19, 12).Value = 0.9932640239640975
$ws.Cells.Item(19, 13).Value = 1.042784359460434
$ws.Cells.Item(19, 14).Value = 1.032216827587152
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.024271366373161
$ws.Cells.Item(20, 4).Value = 1.032684246953516
$ws.Cells.Item(20, 5).Value = 0.9894336180360679
$ws.Cells.Item(20, 6).Value = 1.03916748954583
$ws.Cells.Item(20, 9).Value = 1.031886297319587
$ws.Cells.Item(20, 10).Value = 1.03054749833575
$ws.Cells.Item(20, 11).Value = 1.036069411990714
$ws.Cells.Item(20, 12).Value = 0.9929783193494215
$ws.Cells.Item(20, 13).Value = 1.042530001582685
$ws.Cells.Item(20, 14).Value = 1.032010993988837
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.023199124240644
$ws.Cells.Item(21, 4).Value = 1.031860994920946
$ws.Cells.Item(21, 5).Value = 0.9882828385668249
$ws.Cells.Item(21, 6).Value = 1.038130175167754
$ws.Cells.Item(21, 9).Value = 1.031681592061513
$ws.Cells.Item(21, 10).Value = 1.029878443336643
$ws.Cells.Item(21, 11).Value = 1.035456605191188
$ws.Cells.Item(21, 12).Value = 0.9920501090198102
$ws.Cells.Item(21, 13).Value = 1.041702403623842
$ws.Cells.Item(21, 14).Value = 1.03134098885489
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.022524966737045
$ws.Cells.Item(22, 4).Value = 1.031343014098917
$ws.Cells.Item(22, 5).Value = 0.9875604150241495
$ws.Cells.Item(22, 6).Value = 1.037478147094519
$ws.Cells.Item(22, 9).Value = 1.031550840811468
$ws.Cells.Item(22, 10).Value = 1.029457157834987
$ws.Cells.Item(22, 11).Value = 1.035070221412779
$ws.Cells.Item(22, 12).Value = 0.9914670000341481
$ws.Cells.Item(22, 13).Value = 1.041181562275112
$ws.Cells.Item(22, 14).Value = 1.030919105079499
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.022882349852725
$ws.Cells.Item(23, 4).Value = 1.031617640168086
$ws.Cells.Item(23, 5).Value = 0.9879432794643023
$ws.Cells.Item(23, 6).Value = 1.037823782776921
$ws.Cells.Item(23, 9).Value = 1.031620348159085
$ws.Cells.Item(23, 10).Value = 1.029680548164192
$ws.Cells.Item(23, 11).Value = 1.035275153581206
$ws.Cells.Item(23, 12).Value = 0.991776070289318
$ws.Cells.Item(23, 13).Value = 1.041457716973156
$ws.Cells.Item(23, 14).Value = 1.031142812648595
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.024289342677261
$ws.Cells.Item(24, 4).Value = 1.032698042529532
$ws.Cells.Item(24, 5).Value = 0.9894529299347244
$ws.Cells.Item(24, 6).Value = 1.03918488310389
$ws.Cells.Item(24, 9).Value = 1.031889694520319
$ws.Cells.Item(24, 10).Value = 1.030558704499224
$ws.Cells.Item(24, 11).Value = 1.036079667271709
$ws.Cells.Item(24, 12).Value = 0.9929938892766442
$ws.Cells.Item(24, 13).Value = 1.042543867863303
$ws.Cells.Item(24, 14).Value = 1.032022216066349
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.025922076922809
$ws.Cells.Item(25, 4).Value = 1.033950139586246
$ws.Cells.Item(25, 5).Value = 0.9912096547607049
$ws.Cells.Item(25, 6).Value = 1.040765078415478
$ws.Cells.Item(25, 9).Value = 1.032193331384202
$ws.Cells.Item(25, 10).Value = 1.031575014195638
$ws.Cells.Item(25, 11).Value = 1.037008493826416
$ws.Cells.Item(25, 12).Value = 0.9944092447426414
$ws.Cells.Item(25, 13).Value = 1.043802085190009
$ws.Cells.Item(25, 14).Value = 1.033039969039104
